# Update the cryptos price/volume table with the latest scraped values.
# Cells in column D that look like plain decimal numbers (single dot, e.g.
# "267.06") are prefixed with a leading apostrophe so Excel stores them as
# text (matching the original inlineStr cells) instead of auto-converting
# them to numbers and silently dropping significant trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "19.721.70"
$ws.Range("E2").Value = "  -8.90%  "
$ws.Range("D3").Value = "1.385.52"
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").Value = "'267.06"
$ws.Range("E6").Value = "  -7.40%  "
$ws.Range("D7").Value = "'0.3637"
$ws.Range("E7").Value = "  -7.44%  "
$ws.Range("D8").Value = "'0.3030"
$ws.Range("E8").Value = "  -4.47%  "
$ws.Range("D9").Value = "'38.14"
$ws.Range("E9").Value = "  -10.09%  "
$ws.Range("D10").Value = "'0.9725"
$ws.Range("E10").Value = "  -7.87%  "
$ws.Range("D11").Value = "'0.06397"
$ws.Range("E11").Value = "  -10.86%  "
$ws.Range("D12").Value = "'1.003"
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("D13").Value = "'5.276"
$ws.Range("E13").Value = "  -7.01%  "
$ws.Range("D14").Value = "'6.043"
$ws.Range("E14").Value = "  -8.21%  "
$ws.Range("D15").Value = "'16.43"
$ws.Range("E15").Value = "  -11.44%  "
$ws.Range("D16").Value = "1.390.40"
$ws.Range("E16").Value = "  -10.25%  "
$ws.Range("D17").Value = "'0.000009862"
$ws.Range("E17").Value = "  -9.80%  "
$ws.Range("D18").Value = "'0.05617"
$ws.Range("E18").Value = "  -14.72%  "
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("D20").Value = "'69.60"
$ws.Range("E20").Value = "  -16.85%  "
$ws.Range("D21").Value = "'5.494"
$ws.Range("E21").Value = "  -10.23%  "
$ws.Range("D22").Value = "'14.28"
$ws.Range("E22").Value = "  -7.62%  "
$ws.Range("E23").Value = "  -2.04%  "
$ws.Range("D24").Value = "'2.240"
$ws.Range("E24").Value = "  -4.90%  "
$ws.Range("D25").Value = "19.721.63"
$ws.Range("E25").Value = "  -8.94%  "
$ws.Range("D26").Value = "'2.151"
$ws.Range("E26").Value = "  -8.62%  "
$ws.Range("D27").Value = "'135.47"
$ws.Range("E27").Value = "  -9.31%  "
$ws.Range("D28").Value = "'16.47"
$ws.Range("E28").Value = "  -10.06%  "
$ws.Range("D29").Value = "1.544.10"
$ws.Range("E29").Value = "  -9.83%  "
$ws.Range("D30").Value = "'107.14"
$ws.Range("E30").Value = "  -8.48%  "
$ws.Range("D31").Value = "'3.828"
$ws.Range("E31").Value = "  -21.08%  "
$ws.Range("D32").Value = "'5.189"
$ws.Range("E32").Value = "  -14.57%  "
$ws.Range("D33").Value = "'0.7889"
$ws.Range("E33").Value = "  -15.98%  "
$ws.Range("D34").Value = "'0.07565"
$ws.Range("E34").Value = "  -7.08%  "
$ws.Range("D35").Value = "'8.135"
$ws.Range("E35").Value = "  -4.50%  "
$ws.Range("D36").Value = "'1.002"
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("D37").Value = "'0.05603"
$ws.Range("E37").Value = "  -6.92%  "
$ws.Range("D38").Value = "'4.655"
$ws.Range("E38").Value = "  -9.94%  "
$ws.Range("D39").Value = "'0.02015"
$ws.Range("E39").Value = "  -9.21%  "
$ws.Range("D40").Value = "'0.1862"
$ws.Range("E40").Value = "  -7.98%  "
$ws.Range("D41").Value = "'9.875"
$ws.Range("E41").Value = "  -9.55%  "
$ws.Range("D42").Value = "'1.283"
$ws.Range("E42").Value = "  -11.90%  "
$ws.Range("E43").Value = "  -11.60%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.5157"
$ws.Range("E44").Value = "  -10.57%  "
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").Value = "'3.449"
$ws.Range("E45").Value = "  -6.96%  "
$ws.Range("D46").Value = "'11.89"
$ws.Range("E46").Value = "  -9.00%  "
$ws.Range("D47").Value = "'0.4954"
$ws.Range("E47").Value = "  -10.02%  "
$ws.Range("D48").Value = "'108.34"
$ws.Range("E48").Value = "  -6.62%  "
$ws.Range("D49").Value = "'1.716"
$ws.Range("E49").Value = "  -8.72%  "
$ws.Range("D50").Value = "'1.004"
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("E51").Value = "  -11.67%  "
